# "update after 2nd teaching"
$p = $ppt.ActivePresentation

# --- Slide 1: Subtitle 2 shape - split last line into "  --User login" + "/logout" ---
$s1 = $p.Slides.Item(1)
$sub = $s1.Shapes.Item("Subtitle 2")
$tr1 = $sub.TextFrame.TextRange
$lastPara = $tr1.Paragraphs(2, 1)
# First clear to unrelated text so the final assignment below lands as a single run
# (avoids the engine preserving a stray "common prefix" run from the old text).
$lastPara.Text = "zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"
$lastPara.Text = "  --User login/logout"
# Now split "/logout" into its own run.
$null = $tr1.Replace("/logout", "/logout")

# --- Slide 2: Content Placeholder 2 - "when login" -> "After login" ---
$s2 = $p.Slides.Item(2)
$content = $s2.Shapes.Item("Content Placeholder 2")
$ctr = $content.TextFrame.TextRange
$null = $ctr.Replace("when login, all activity assume the ", "After login, all activity assume the ")

# --- Slide 11: Straight Connector 6 - move line ---
$s11 = $p.Slides.Item(11)
$conn = $s11.Shapes.Item("Straight Connector 6")
$conn.Left = 2127123 / 12700
$conn.Top = 2121409 / 12700
